$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-26 14:57:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-26 14:57:27"
$wsZhCn.Range("K4").Value = "2016-08-26 14:57:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-26 14:57:32"
$wsDeDe.Range("K4").Value = "2016-08-26 14:57:52"
